$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-18 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-19 Friday", 2) | Out-Null
$d.Content.Find.Execute("16×33=528", $true, $false, $false, $false, $false, $true, 1, $false, "68×13=884", 2) | Out-Null
$d.Content.Find.Execute("99×49=4851", $true, $false, $false, $false, $false, $true, 1, $false, "49×14=686", 2) | Out-Null
$d.Content.Find.Execute("55×49=2695", $true, $false, $false, $false, $false, $true, 1, $false, "14×74=1036", 2) | Out-Null
$d.Content.Find.Execute("45×47=2115", $true, $false, $false, $false, $false, $true, 1, $false, "78×43=3354", 2) | Out-Null
$d.Content.Find.Execute("95×46=4370", $true, $false, $false, $false, $false, $true, 1, $false, "50×31=1550", 2) | Out-Null
$d.Content.Find.Execute("99×38=3762", $true, $false, $false, $false, $false, $true, 1, $false, "68×59=4012", 2) | Out-Null
$d.Content.Find.Execute("92×47=4324", $true, $false, $false, $false, $false, $true, 1, $false, "29×37=1073", 2) | Out-Null
$d.Content.Find.Execute("13×92=1196", $true, $false, $false, $false, $false, $true, 1, $false, "31×11=341", 2) | Out-Null
$d.Content.Find.Execute("14×60=840", $true, $false, $false, $false, $false, $true, 1, $false, "51×48=2448", 2) | Out-Null
$d.Content.Find.Execute("65×55=3575", $true, $false, $false, $false, $false, $true, 1, $false, "14×17=238", 2) | Out-Null
$d.Content.Find.Execute("65×22=1430", $true, $false, $false, $false, $false, $true, 1, $false, "67×30=2010", 2) | Out-Null
$d.Content.Find.Execute("35×28=980", $true, $false, $false, $false, $false, $true, 1, $false, "47×92=4324", 2) | Out-Null
$d.Content.Find.Execute("92×96=8832", $true, $false, $false, $false, $false, $true, 1, $false, "81×92=7452", 2) | Out-Null
$d.Content.Find.Execute("73×35=2555", $true, $false, $false, $false, $false, $true, 1, $false, "12×50=600", 2) | Out-Null
$d.Content.Find.Execute("75×51=3825", $true, $false, $false, $false, $false, $true, 1, $false, "50×32=1600", 2) | Out-Null
$d.Content.Find.Execute("77×27=2079", $true, $false, $false, $false, $false, $true, 1, $false, "59×81=4779", 2) | Out-Null
$d.Content.Find.Execute("79×20=1580", $true, $false, $false, $false, $false, $true, 1, $false, "46×87=4002", 2) | Out-Null
$d.Content.Find.Execute("67×23=1541", $true, $false, $false, $false, $false, $true, 1, $false, "71×91=6461", 2) | Out-Null
$d.Content.Find.Execute("26×32=832", $true, $false, $false, $false, $false, $true, 1, $false, "40×36=1440", 2) | Out-Null
$d.Content.Find.Execute("31×64=1984", $true, $false, $false, $false, $false, $true, 1, $false, "13×18=234", 2) | Out-Null
$d.Content.Find.Execute("82×36=2952", $true, $false, $false, $false, $false, $true, 1, $false, "77×30=2310", 2) | Out-Null
$d.Content.Find.Execute("76×63=4788", $true, $false, $false, $false, $false, $true, 1, $false, "74×13=962", 2) | Out-Null
$d.Content.Find.Execute("91×43=3913", $true, $false, $false, $false, $false, $true, 1, $false, "88×89=7832", 2) | Out-Null
$d.Content.Find.Execute("82×51=4182", $true, $false, $false, $false, $false, $true, 1, $false, "93×76=7068", 2) | Out-Null
$d.Content.Find.Execute("81×89=7209", $true, $false, $false, $false, $false, $true, 1, $false, "39×90=3510", 2) | Out-Null
